# Corrections in "Perfil de Empleados" test sheet:
# Add a new "OK" column (D) marking the test cases that passed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "OK"
$ws.Range("D5").Value = "OK"
$ws.Range("D7").Value = "OK"
$ws.Range("D8").Value = "OK"

$ws.Range("D6").Select() | Out-Null
